# LinkedIn carousel -> single-slide draft + prompts
# Applies:
#   1. Delete slides 2-6 (keep only slide 1)
#   2. Update slide 1 title text + size
#   3. Update slide 1 body bullets (swap in the Google-News blurb, add 2 follow-up prompts)
#   4. Reposition/resize the slide 1 picture

$p = $ppt.ActivePresentation

# 1) Drop slides 2..6, leaving a single-slide deck.
for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}

$s = $p.Slides.Item(1)

# 2) Title: new headline, bumped up to 32pt (bold already inherited from the
#    placeholder's default run properties, so it only needs the size bump).
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = 'NTPC Green Energy board okays 50:50 JV with GAIL'
$titleRange.Font.Size = 32

# 3) Body placeholder: replace the two existing bullet paragraphs and append a third.
$body = $s.Shapes.Item(2)
$bodyRange = $body.TextFrame.TextRange

$bodyRange.Paragraphs(2, 1).Text = '<a href="https://news.google.com/rss/articles/CBMilwFBVV95cUxOTy1QWnlsUWdnZzIxMS1PN01PRjA2ZFlUbEV0eHYwODRncm5VekVBS214TDUwMk9GLWZlVGZ6NUlZMFd1T3BlcGVoY0FTbk9selhyc21QSkdiekl6VkRXb3lXT2tNUGpQWXN5dFBTQzJSZHhHYUI1OV9ZdTNsMWk2dnY5VmRZZElPbGpBSHJNY3lSWU5JcHFN?oc=5" target="_blank">NTPC Green Energy board okays 50:50 JV with GAIL</a>&nbsp;&nbsp;<font color="#6f6f6f">India Infoline</font>'

$bodyRange.Paragraphs(3, 1).Text = 'Strategic boost to India’s renewable ecosystem'

$null = $bodyRange.InsertAfter("`rSupports long-term clean energy transition")

# 4) Picture: recenter/resize on the slide.
$pic = $s.Shapes.Item(3)
$pic.Left = 72
$pic.Top = 144
$pic.Width = 576
$pic.Height = 252
